$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add new column D: "PrimeNumbersOptimal" results (Time + Instr) ---

# Header style: copy the big bold/centered header format from the existing
# "PrimeNumbersPositive - ..." header cells (C1 / C13) onto D1 / D13.
$ws.Range("C1").Copy() | Out-Null
$ws.Range("D1").PasteSpecial(-4122) | Out-Null

$ws.Range("C13").Copy() | Out-Null
$ws.Range("D13").PasteSpecial(-4122) | Out-Null

# Data style for the "Time" block (rows 2-9): matches the plain-number
# format used in column E of that block (no special alignment).
$ws.Range("E3").Copy() | Out-Null
$ws.Range("D2:D9").PasteSpecial(-4122) | Out-Null

# Data style for the "Instr" block (rows 14-21): matches the right-aligned
# integer format already used by columns A-C there.
$ws.Range("C14").Copy() | Out-Null
$ws.Range("D14:D21").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = $false

# Header text
$ws.Range("D1").Value = "PrimeNumbersOptimal - Time"
$ws.Range("D13").Value = "PrimeNumbersOptimal - Instr"

# "Time" measurements (rows 2-9)
$ws.Range("D2").Value = 0
$ws.Range("D3").Value = 0
$ws.Range("D4").Value = 0
$ws.Range("D5").Value = 0
$ws.Range("D6").Value = 0
$ws.Range("D7").Value = 1
$ws.Range("D8").Value = 5
$ws.Range("D9").Value = 15

# "Instr" (instruction count) measurements (rows 14-21)
$ws.Range("D14").Value = 108
$ws.Range("D15").Value = 337
$ws.Range("D16").Value = 1061
$ws.Range("D17").Value = 3351
$ws.Range("D18").Value = 10591
$ws.Range("D19").Value = 33487
$ws.Range("D20").Value = 105892
$ws.Range("D21").Value = 334855

# Widen column D so the new header text fits (bestFit-style width, close to
# the author's 39.140625 chars).
$ws.Columns("D").ColumnWidth = 38.25

# Move the selection, as recorded in the saved workbook.
$ws.Range("D10").Select() | Out-Null
